# "Start using the status of a graphic"
# Fill in the two blank timesheet rows (56 & 57) on the "2018" sheet with
# the new tasks, and give row 57 a manual formula for the hosting renewal
# charge instead of the usual Hours*Rate calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018")

# Row 56 - Added Miro Prints plus updates
$ws.Range("A56").Value = "Added Miro Prints plus updates"
$ws.Range("B56").Value = 43734
$ws.Range("C56").Value = 0.5
$ws.Range("D56").Value = 25

# Row 57 - Renewal of Website Hosting
$ws.Range("A57").Value = "Renewal of Website Hosting"
$ws.Range("B57").Value = 43734
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 0
$ws.Range("E57").Formula = "=(59.4+19.95)*0.81"
